$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.348.06"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "3.156.91"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.79"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.10"
$ws.Range("E6").Value = "  -5.96%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.157.10"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -3.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  -6.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.54"
$ws.Range("E11").Value = "  -4.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.476"
$ws.Range("E12").Value = "  -4.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  -5.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.24"
$ws.Range("E14").Value = "  -7.00%  "
$ws.Range("D15").Value = "3.669.41"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").Value = "64.309.00"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "3.155.17"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.98"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.64"
$ws.Range("E20").Value = "  -5.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.64"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.77"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.84"
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.51"
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.46"
$ws.Range("E28").Value = "  -6.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -6.39%  "
$ws.Range("E30").Value = "  -30.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.83"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("E32").Value = "  -4.91%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.31"
$ws.Range("E34").Value = "  -6.74%  "
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.04"
$ws.Range("E36").Value = "  -4.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.21"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("D38").Value = "0.0₃0732"
$ws.Range("E38").Value = "  -7.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "454.20"
$ws.Range("E39").Value = "  -8.07%  "
$ws.Range("E40").Value = "  -6.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0398"
$ws.Range("E41").Value = "  -5.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.120"
$ws.Range("E42").Value = "  -6.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.48"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("D44").Value = "2.858.11"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.270"
$ws.Range("E45").Value = "  -7.76%  "
$ws.Range("E46").Value = "  -7.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.48"
$ws.Range("E47").Value = "  -6.12%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.32"
$ws.Range("E49").Value = "  -4.10%  "
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.72"
$ws.Range("E51").Value = "  -1.58%  "
